# Apply the edit described by the diff:
# Insert a new row at row 143 (pushing old rows 143..224 down to 144..225)
# and populate the new row 143 with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 143. Excel will copy the
# formatting (including the date NumberFormat on column D) from the
# row above, matching the original sheet's layout.
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with the new record's data.
$ws.Cells.Item(143, 1).Value = 11
$ws.Cells.Item(143, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(143, 3).Value = "Bíobío"
$ws.Cells.Item(143, 4).Value = 45141
$ws.Cells.Item(143, 5).Value = 8
$ws.Cells.Item(143, 6).Value = 100112021
$ws.Cells.Item(143, 7).Value = "Ají"
$ws.Cells.Item(143, 8).Value = "Inferno"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 50
$ws.Cells.Item(143, 11).Value = 17000
$ws.Cells.Item(143, 12).Value = 18000
$ws.Cells.Item(143, 13).Value = 17400
$ws.Cells.Item(143, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(143, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(143, 16).Value = 1740
$ws.Cells.Item(143, 17).Value = 10
$ws.Cells.Item(143, 18).Value = "Hortaliza"
